# The deck's single slide master currently carries the "Integral" / "Red
# Violet" theme (ppt/theme/theme2.xml) as its applied color scheme, while an
# unused "Office Theme" color scheme sits alongside it (ppt/theme/theme1.xml,
# only wired to the notes master). The edit re-applies the built-in default
# Office color scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) to the
# presentation's theme, replacing the Red Violet palette with the standard
# Office one.

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.ColorScheme

# Office (default) theme color scheme, in clrScheme document order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$scheme.Colors(1).RGB  = 0          # dk1      000000
$scheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2      44546A
$scheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$scheme.Colors(10).RGB = 4697456    # accent6  70AD47
$scheme.Colors(11).RGB = 12673797   # hlink    0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink 954F72
